$d = $word.ActiveDocument

# 1. Replace the placeholder ID text, and remove the trailing space run
$d.Content.Find.Execute("**ID__AFFARS_5311_topic_2__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5311_1__ID**", 2)

# 2. Update the indentation and add a paragraph border (space only, no line) on the first paragraph
$p1 = $d.Paragraphs(1)
$p1.Format.LeftIndent = 11.25

$b = $p1.Format.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5
